$wb = $excel.ActiveWorkbook

# Sheet 1 (ALC)
$ws = $wb.Worksheets.Item(1)
$ws.Range("H53").Value = 100697.75
$ws.Range("J53").Value = 154305.92
$ws.Range("L53").Value = 154305.92
$ws.Range("N53").Value = -155579.92
$ws.Range("H70").Value = 4249.75
$ws.Range("J70").Value = 4333.3335
$ws.Range("L70").Value = 13000.0005
$ws.Range("N70").Value = -13540.0005
$ws.Range("H73").Value = 4249.75
$ws.Range("J73").Value = 4333.3335
$ws.Range("L73").Value = 13000.0005
$ws.Range("N73").Value = -14872.0005
$ws.Range("H98").Value = 614.6667
$ws.Range("I98").Value = 506.7647
$ws.Range("K98").Value = 506.7647
$ws.Range("M98").Value = 991.2353000000001
$ws.Range("H122").Value = 614.6667
$ws.Range("I122").Value = 506.7647
$ws.Range("K122").Value = 1520.2941
$ws.Range("M122").Value = 929.7058999999999
$ws.Range("H127").Value = 1194.5
$ws.Range("I127").Value = 592.6667
$ws.Range("K127").Value = 1778.0001
$ws.Range("M127").Value = 3181.9999
$ws.Range("H132").Value = 17029.188
$ws.Range("I132").Value = 16955.158
$ws.Range("K132").Value = 50865.474
$ws.Range("M132").Value = -48335.474
$ws.Range("H137").Value = 4291.3794
$ws.Range("I137").Value = 4290.385
$ws.Range("J137").Value = 4300
$ws.Range("K137").Value = 12871.155
$ws.Range("L137").Value = 12900
$ws.Range("M137").Value = -10321.155
$ws.Range("N137").Value = -18000
$ws.Range("H138").Value = 2799.1296
$ws.Range("J138").Value = 4395.56
$ws.Range("L138").Value = 13186.68
$ws.Range("N138").Value = -23466.68

# Sheet 2 (ARM)
$ws = $wb.Worksheets.Item(2)
$ws.Range("H102").Value = 3473.6
$ws.Range("I102").Value = 2548.1428
$ws.Range("K102").Value = 2548.1428
$ws.Range("M102").Value = -926.1428000000001
$ws.Range("H110").Value = 4175.3438
$ws.Range("I110").Value = 3346.8215
$ws.Range("K110").Value = 3346.8215
$ws.Range("M110").Value = -1301.8215
$ws.Range("H122").Value = 1818.3
$ws.Range("I122").Value = 1476.5333
$ws.Range("K122").Value = 4429.5999
$ws.Range("M122").Value = -1979.5999

# Sheet 3 (BSM)
$ws = $wb.Worksheets.Item(3)
$ws.Range("H23").Value = 1014
$ws.Range("J23").Value = 1014
$ws.Range("L23").Value = 1014
$ws.Range("N23").Value = -1580
$ws.Range("H99").Value = 4319.7617
$ws.Range("I99").Value = 3670.6667
$ws.Range("K99").Value = 3670.6667
$ws.Range("M99").Value = -2172.6667
$ws.Range("H134").Value = 6224.909
$ws.Range("I134").Value = 6045.143
$ws.Range("K134").Value = 18135.429
$ws.Range("M134").Value = -15600.429

# Sheet 4 (CRP)
$ws = $wb.Worksheets.Item(4)
$ws.Range("H16").Value = 2945.6897
$ws.Range("J16").Value = 4861.778
$ws.Range("L16").Value = 4861.778
$ws.Range("N16").Value = -5435.778
$ws.Range("H31").Value = 5452.8
$ws.Range("I31").Value = 2562.7856
$ws.Range("J31").Value = 9131
$ws.Range("K31").Value = 2562.7856
$ws.Range("L31").Value = 9131
$ws.Range("M31").Value = -2267.7856
$ws.Range("N31").Value = -9721
$ws.Range("H33").Value = 1000
$ws.Range("I33").Value = 1000
$ws.Range("K33").Value = 1000
$ws.Range("M33").Value = -621
$ws.Range("H34").Value = 5452.8
$ws.Range("I34").Value = 2562.7856
$ws.Range("J34").Value = 9131
$ws.Range("K34").Value = 2562.7856
$ws.Range("L34").Value = 9131
$ws.Range("M34").Value = -2360.7856
$ws.Range("N34").Value = -9535
$ws.Range("H41").Value = 20872.5
$ws.Range("J41").Value = 24496.666
$ws.Range("L41").Value = 24496.666
$ws.Range("N41").Value = -25352.666
$ws.Range("H113").Value = 2945.6897
$ws.Range("J113").Value = 4861.778
$ws.Range("L113").Value = 4861.778
$ws.Range("N113").Value = -9201.778
$ws.Range("H132").Value = 7506.839
$ws.Range("I132").Value = 7060.852
$ws.Range("K132").Value = 21182.556
$ws.Range("M132").Value = -18652.556
$ws.Range("H134").Value = 3734.9678
$ws.Range("I134").Value = 3138.2964
$ws.Range("K134").Value = 9414.889200000001
$ws.Range("M134").Value = -6879.889200000001

# Sheet 5 (CUL)
$ws = $wb.Worksheets.Item(5)
$ws.Range("H47").Value = 3930.2
$ws.Range("I47").Value = 3940.25
$ws.Range("K47").Value = 11820.75
$ws.Range("M47").Value = -11389.75
$ws.Range("H75").Value = 225
$ws.Range("J75").Value = 225
$ws.Range("L75").Value = 675
$ws.Range("N75").Value = -2671
$ws.Range("H78").Value = 225
$ws.Range("J78").Value = 225
$ws.Range("L78").Value = 2025
$ws.Range("N78").Value = -12009
$ws.Range("H92").Value = 257.25
$ws.Range("J92").Value = 243
$ws.Range("L92").Value = 729
$ws.Range("N92").Value = -3225
$ws.Range("H131").Value = 13890983
$ws.Range("I131").Value = 71429544
$ws.Range("K131").Value = 214288632
$ws.Range("M131").Value = -214283592

# Sheet 6 (GSM)
$ws = $wb.Worksheets.Item(6)
$ws.Range("H31").Value = 3499.3
$ws.Range("I31").Value = 2999.2222
$ws.Range("K31").Value = 2999.2222
$ws.Range("M31").Value = -2707.2222
$ws.Range("H37").Value = 3499.3
$ws.Range("I37").Value = 2999.2222
$ws.Range("K37").Value = 2999.2222
$ws.Range("M37").Value = -2722.2222
$ws.Range("H126").Value = 3929
$ws.Range("I126").Value = 2935.1667
$ws.Range("K126").Value = 8805.500100000001
$ws.Range("M126").Value = -6335.500100000001
$ws.Range("H132").Value = 1923.9166
$ws.Range("J132").Value = 1868.1
$ws.Range("L132").Value = 5604.299999999999
$ws.Range("N132").Value = -10664.3
$ws.Range("H133").Value = 101000
$ws.Range("I133").Value = 101000
$ws.Range("K133").Value = 101000
$ws.Range("M133").Value = -95940

# Sheet 7 (LTW)
$ws = $wb.Worksheets.Item(7)
$ws.Range("H22").Value = 1484.4445
$ws.Range("J22").Value = 641
$ws.Range("L22").Value = 641
$ws.Range("N22").Value = -1231
$ws.Range("H26").Value = 60000
$ws.Range("I26").Value = 60000
$ws.Range("K26").Value = 60000
$ws.Range("M26").Value = -59705
$ws.Range("H27").Value = 1484.4445
$ws.Range("J27").Value = 641
$ws.Range("L27").Value = 641
$ws.Range("N27").Value = -855
$ws.Range("H40").Value = 3898.25
$ws.Range("I40").Value = 3898.25
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3898.25
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -3762.25
$ws.Range("H132").Value = 16009.127
$ws.Range("I132").Value = 13846.552
$ws.Range("J132").Value = 25657.54
$ws.Range("K132").Value = 41539.656
$ws.Range("L132").Value = 76972.62
$ws.Range("M132").Value = -39009.656
$ws.Range("N132").Value = -82032.62
$ws.Range("N40").ClearContents()

# Sheet 8 (WVR)
$ws = $wb.Worksheets.Item(8)
$ws.Range("H62").Value = 8873.454
$ws.Range("I62").Value = 8567
$ws.Range("K62").Value = 8567
$ws.Range("M62").Value = -7943
$ws.Range("H65").Value = 8873.454
$ws.Range("I65").Value = 8567
$ws.Range("K65").Value = 42835
$ws.Range("M65").Value = -39715
$ws.Range("H122").Value = 4561
$ws.Range("I122").Value = 3526.3333
$ws.Range("K122").Value = 10578.9999
$ws.Range("M122").Value = -8128.999899999999
